$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# A1 holds a date serial number; bump it by one day (45310 -> 45311)
$ws.Range("A1").Value = 45311

# Update prices in column D for rows 33-35
$ws.Range("D33").Value = 161.5
$ws.Range("D34").Value = 190.8
$ws.Range("D35").Value = 211
